$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# 1. Insert a brand-new empty paragraph right after paragraph 2
#    ("Geographic Location Attribute Predictor System").
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$r2 = $p2.Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()

# The freshly minted paragraph is now paragraph 3; make sure it is a truly
# empty <w:p/> (no stray run) by round-tripping it through InsertXML.
$p3 = $d.Paragraphs.Item(3)
$p3.Range.InsertXML("<w:p $wns/>")

# ---------------------------------------------------------------------------
# 2. Paragraph 4 ("We have created G.L.A.P.S ...") becomes the new
#    multi-run "G.L.A.P.S will be a web application ..." paragraph.
# ---------------------------------------------------------------------------
$para4xml = '<w:p ' + $wns + '>' +
  '<w:r><w:t>G.L.A.P.S</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> will be a web application that will predict the future value of homes in a </w:t></w:r>' +
  '<w:proofErr w:type="gramStart"/>' +
  '<w:r><w:t>particular area</w:t></w:r>' +
  '<w:proofErr w:type="gramEnd"/>' +
  '<w:r><w:t xml:space="preserve">. </w:t></w:r>' +
  '<w:r><w:t>T</w:t></w:r>' +
  '<w:r><w:t>he software</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:t>use</w:t></w:r>' +
  '<w:r><w:t>s</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> data from the United States Census to predict the future values of homes in the Fayetteville area. A minor league</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> baseball</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> stadium is scheduled to open in Cumberland County (Fayetteville) in April of 2019. We have gathered data from counties in other areas of the U.S. that also had minor league </w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">baseball </w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">stadiums open within the past </w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">eight years. Our data will be held in a sqlite3 database. </w:t></w:r>' +
  '</w:p>'

$p4 = $d.Paragraphs.Item(4)
$p4.Range.InsertXML($para4xml)

# ---------------------------------------------------------------------------
# 3. Paragraph 5 ("For this version of our software ...") becomes the new
#    "We have gathered data from three years before the stadium..." text,
#    keeping the _GoBack bookmark in place.
# ---------------------------------------------------------------------------
$para5xml = '<w:p ' + $wns + '>' +
  '<w:r><w:t>We have gathered data from three years before the stadium</w:t></w:r>' +
  '<w:r><w:t>s</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
  '<w:bookmarkEnd w:id="0"/>' +
  '<w:r><w:t xml:space="preserve"> opened and three years after. We will be feeding this data into to a Python machine learning program to forecast the values of homes in Fayetteville three years after the stadium has opened its doors. </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>Tensorflow</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> and </w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:t>Keras</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '<w:r><w:t xml:space="preserve"> are the Python machine learning libraries we will use. We plan to try them both and then utilize whichever is more efficient. We may end up using only one or a combination of both.</w:t></w:r>' +
  '</w:p>'

$p5 = $d.Paragraphs.Item(5)
$p5.Range.InsertXML($para5xml)

# ---------------------------------------------------------------------------
# 4. Two brand-new paragraphs get added after paragraph 5, re-using the
#    first of the three trailing empty paragraphs (paragraph 6).
# ---------------------------------------------------------------------------
$para6xml = '<w:p ' + $wns + '><w:r><w:t>G.L.A.P.S will be helpful for those that are interested in buying or selling a house in a certain area, those that are deciding where to move to or those that wish to predict the future value of a home that they own. It may also be beneficial to local government entities to predict the growth of an area or tax values.</w:t></w:r></w:p>'

$para7xml = '<w:p ' + $wns + '><w:r><w:t>The user will enter an address into the web-based front end and a value will be returned to the user. Initially, we will only predict the value of homes and we will stick to areas that relate to minor league baseball stadiums. In the future, we hope to broaden the results of G.L.A.P.S by allowing the user to search for a variety of attributes such as predicted changes in population, income levels, education levels and property tax rates. We would also like to expand the system so that we are not limited to only areas related to baseball stadiums.</w:t></w:r></w:p>'

$p6 = $d.Paragraphs.Item(6)
$p6.Range.InsertXML($para6xml)

$p6b = $d.Paragraphs.Item(6)
$r6b = $p6b.Range
$r6b.Collapse(0)
$r6b.InsertParagraphAfter()

$p7 = $d.Paragraphs.Item(7)
$p7.Range.InsertXML($para7xml)
